$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 have their data (date, volume, prices) swapped.

# D: Fecha
$ws.Range("D2").Value = 44421
$ws.Range("D3").Value = 44291

# M: Volumen
$ws.Range("M2").Value = 30
$ws.Range("M3").Value = 15

# N: Precio minimo
$ws.Range("N2").Value = 24000
$ws.Range("N3").Value = 23000

# O: Precio maximo
$ws.Range("O2").Value = 24000
$ws.Range("O3").Value = 23000

# P: Precio promedio ponderado
$ws.Range("P2").Value = 24000
$ws.Range("P3").Value = 23000

# S: Precio $/Kg
$ws.Range("S2").Value = 1200
$ws.Range("S3").Value = 1150
